$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text storage (avoids Excel
# auto-converting numeric-looking strings like "1.001" or "244.09" into
# real numbers), and restoring the original cell style afterwards so no
# stray formatting (e.g. quote-prefix/text number format) is left behind.
function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "29.148.98"
$ws.Range("E2").Value = "  +0.23%  "
Set-TextValue $ws.Range("D3") "1.840.67"
$ws.Range("E3").Value = "  +0.27%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "244.09"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("E6").Value = "  -1.19%  "
Set-TextValue $ws.Range("D7") "1.002"
$ws.Range("E7").Value = "  +0.09%  "
Set-TextValue $ws.Range("D8") "0.07510"
$ws.Range("E8").Value = "  -0.58%  "
Set-TextValue $ws.Range("D9") "0.2939"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  +1.47%  "
Set-TextValue $ws.Range("D11") "0.07717"
$ws.Range("E11").Value = "  -0.38%  "
Set-TextValue $ws.Range("D12") "1.880.74"
$ws.Range("E12").Value = "  +2.47%  "
Set-TextValue $ws.Range("D13") "5.023"
$ws.Range("E13").Value = "  +0.30%  "
Set-TextValue $ws.Range("D14") "0.6763"
$ws.Range("E14").Value = "  +0.65%  "
Set-TextValue $ws.Range("D15") "83.07"
$ws.Range("E15").Value = "  -0.33%  "
Set-TextValue $ws.Range("D16") "0.000009267"
$ws.Range("E16").Value = "  -3.94%  "
Set-TextValue $ws.Range("D17") "5.980"
$ws.Range("E17").Value = "  -1.79%  "
Set-TextValue $ws.Range("D18") "29.152.42"
$ws.Range("E18").Value = "  +0.16%  "
Set-TextValue $ws.Range("D19") "2.124.71"
$ws.Range("E19").Value = "  +1.81%  "
Set-TextValue $ws.Range("D20") "230.69"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("E21").Value = "  +0.84%  "
Set-TextValue $ws.Range("D22") "1.003"
$ws.Range("E22").Value = "  +0.25%  "
Set-TextValue $ws.Range("D23") "7.189"
$ws.Range("E23").Value = "  -0.30%  "
Set-TextValue $ws.Range("D24") "1.002"
$ws.Range("E24").Value = "  +0.11%  "
Set-TextValue $ws.Range("D25") "160.52"
$ws.Range("E25").Value = "  -0.20%  "
Set-TextValue $ws.Range("D26") "8.562"
$ws.Range("E26").Value = "  +0.07%  "
Set-TextValue $ws.Range("D27") "0.1391"
$ws.Range("E27").Value = "  -0.73%  "
Set-TextValue $ws.Range("D28") "17.91"
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("E29").Value = "  -0.23%  "
Set-TextValue $ws.Range("D30") "4.193"
$ws.Range("E30").Value = "  +1.54%  "
Set-TextValue $ws.Range("D31") "4.148"
$ws.Range("E31").Value = "  +1.66%  "
Set-TextValue $ws.Range("D32") "0.05558"
$ws.Range("E32").Value = "  +3.03%  "
Set-TextValue $ws.Range("D34") "0.7501"
$ws.Range("E34").Value = "  +0.47%  "
Set-TextValue $ws.Range("D35") "1.856"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "2.775"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("B39").Value = "Maker"
$ws.Range("C39").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D39") "1.230.40"
$ws.Range("E39").Value = "  -1.15%  "
$ws.Range("E40").Value = "  +0.09%  "
Set-TextValue $ws.Range("D41") "6.564"
$ws.Range("E41").Value = "  -1.07%  "
Set-TextValue $ws.Range("D42") "0.9021"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  +0.00%  "
Set-TextValue $ws.Range("D44") "2.016.04"
$ws.Range("E44").Value = "  +1.48%  "
Set-TextValue $ws.Range("D45") "102.24"
$ws.Range("E45").Value = "  +0.18%  "
Set-TextValue $ws.Range("D46") "66.41"
$ws.Range("E46").Value = "  +2.23%  "
Set-TextValue $ws.Range("D47") "0.00000000121"
$ws.Range("E47").Value = "  -1.81%  "
Set-TextValue $ws.Range("D48") "0.5104"
$ws.Range("E48").Value = "  -0.29%  "
Set-TextValue $ws.Range("D49") "0.4088"
$ws.Range("E49").Value = "  -0.26%  "
Set-TextValue $ws.Range("D50") "9.131"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +1.04%  "
